# Update column G ("K" - strikeouts) values for rows 2-34 on Sheet1.
# These values were regenerated upstream (save_data regen: use K instead of
# Strike#, regen std/mean, calc and write s_vals); here we just write the
# resulting literal values into the worksheet cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 4
    3  = 4
    4  = 2
    5  = 6
    6  = 2
    7  = 8
    8  = 7
    9  = 9
    10 = 0
    11 = 0
    12 = 4
    13 = 4
    14 = 2
    15 = 9
    16 = 5
    17 = 5
    18 = 4
    19 = 6
    20 = 6
    21 = 6
    22 = 4
    23 = 8
    24 = 5
    25 = 8
    26 = 4
    27 = 7
    28 = 6
    29 = 7
    30 = 3
    31 = 1
    32 = 9
    33 = 7
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
